{"js": "// Load all body paragraphs along with their text so we can locate the\n// paragraphs that need to change by content (more robust than hard-coded\n// indices, in case the document shifts slightly).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  items[i].load(\"text\");\n}\nawait context.sync();\n\n// --- Change 1 -----------------------------------------------------------\n// Collapse the three \"CORE COMPETENCIES\" detail paragraphs into a single\n// condensed summary paragraph.\nlet coreIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Statistical Analysis & Machine Learning: Advanced Statistical Modeling\") === 0) {\n    coreIdx = i;\n    break;\n  }\n}\n\nif (coreIdx !== -1) {\n  const firstPara = items[coreIdx];\n  const secondPara = items[coreIdx + 1];\n  const thirdPara = items[coreIdx + 2];\n\n  // Replace the text of the first paragraph with the condensed summary,\n  // then remove the two paragraphs that followed it.\n  firstPara.insertText(\n    \"Statistical Analysis & Machine Learning \\u2022 Big Data & Data Engineering \\u2022 Data Visualization & Reporting\",\n    Word.InsertLocation.replace\n  );\n  secondPara.delete();\n  thirdPara.delete();\n  await context.sync();\n}\n\n// --- Change 2 -----------------------------------------------------------\n// Insert a new \"TECHNICAL SKILLS\" section (heading + three detail lines)\n// right after the \"Led multi-million dollar research projects...\" bullet.\nlet anchorIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"\\u2022 Led multi-million dollar research projects\") === 0) {\n    anchorIdx = i;\n    break;\n  }\n}\n\nif (anchorIdx !== -1) {\n  const anchor = items[anchorIdx];\n\n  // Insert each new paragraph directly after the anchor; inserting in\n  // reverse order keeps them in the correct final sequence since each one\n  // is placed immediately after the (unmoving) anchor paragraph.\n  const dataViz = anchor.insertParagraph(\n    \"DATA VISUALIZATION & REPORTING Data Visualization; Geospatial Analysis; Interactive Dashboards; Business Intelligence\",\n    Word.InsertLocation.after\n  );\n  const bigData = anchor.insertParagraph(\n    \"BIG DATA & DATA ENGINEERING Big Data Processing; Data Warehousing; Cloud Platforms; Data Pipeline Optimization\",\n    Word.InsertLocation.after\n  );\n  const statAnalysis = anchor.insertParagraph(\n    \"STATISTICAL ANALYSIS & MACHINE LEARNING Advanced Statistical Modeling; Predictive Analytics; Data Mining; Machine Learning\",\n    Word.InsertLocation.after\n  );\n  const heading = anchor.insertParagraph(\"TECHNICAL SKILLS\", Word.InsertLocation.after);\n  heading.style = \"Heading 2\";\n\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1 -------------------------------------------------------------\n# Collapse the three \"CORE COMPETENCIES\" detail paragraphs into a single\n# condensed summary paragraph. Locate the first of the three by its\n# (stable) leading text so this is resilient to small shifts elsewhere in\n# the document.\n$count = $d.Paragraphs.Count\n$coreIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"Statistical Analysis & Machine Learning: Advanced Statistical Modeling*\") {\n        $coreIdx = $i\n        break\n    }\n}\n\nif ($coreIdx -ne -1) {\n    $firstPara = $d.Paragraphs.Item($coreIdx)\n    $firstPara.Range.Text = \"Statistical Analysis & Machine Learning \u2022 Big Data & Data Engineering \u2022 Data Visualization & Reporting\"\n\n    # The two follow-on detail paragraphs now immediately follow the\n    # (just-rewritten) summary paragraph; remove both.\n    $d.Paragraphs.Item($coreIdx + 1).Range.Delete()\n    $d.Paragraphs.Item($coreIdx + 1).Range.Delete()\n}\n\n# --- Change 2 -------------------------------------------------------------\n# Insert a new \"TECHNICAL SKILLS\" section (heading + three detail lines)\n# right after the \"Led multi-million dollar research projects...\" bullet.\n$count = $d.Paragraphs.Count\n$anchorIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Led multi-million dollar research projects*\") {\n        $anchorIdx = $i\n        break\n    }\n}\n\nif ($anchorIdx -ne -1) {\n    $anchor = $d.Paragraphs.Item($anchorIdx)\n    $anchor.Range.InsertParagraphAfter()\n\n    $heading = $d.Paragraphs.Item($anchorIdx + 1)\n    $heading.Range.Text = \"TECHNICAL SKILLS\"\n    $heading.Style = \"Heading 2\"\n\n    $heading.Range.InsertParagraphAfter()\n    $statPara = $d.Paragraphs.Item($anchorIdx + 2)\n    $statPara.Style = \"Normal\"\n    $statPara.Range.Text = \"STATISTICAL ANALYSIS & MACHINE LEARNING Advanced Statistical Modeling; Predictive Analytics; Data Mining; Machine Learning\"\n\n    $statPara.Range.InsertParagraphAfter()\n    $bigDataPara = $d.Paragraphs.Item($anchorIdx + 3)\n    $bigDataPara.Style = \"Normal\"\n    $bigDataPara.Range.Text = \"BIG DATA & DATA ENGINEERING Big Data Processing; Data Warehousing; Cloud Platforms; Data Pipeline Optimization\"\n\n    $bigDataPara.Range.InsertParagraphAfter()\n    $dataVizPara = $d.Paragraphs.Item($anchorIdx + 4)\n    $dataVizPara.Style = \"Normal\"\n    $dataVizPara.Range.Text = \"DATA VISUALIZATION & REPORTING Data Visualization; Geospatial Analysis; Interactive Dashboards; Business Intelligence\"\n}\n"}
